# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    existing "2022-Q2" sheet), by duplicating the "2022-Q2" sheet (so the
#    header row + column-A/header styling come along for free) and then
#    overwriting every data cell with the 2022-Q3 numbers.
# 2. Add a new top data row to "总计" for 2022-Q3 and renumber the
#    existing rows' running index.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q3" sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($null, $summary)

$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q3"

# The source sheet only had 5 rows of fund data (rows 2-5 + header);
# 2022-Q3 needs 6 funds (rows 2-7). Extend by cloning the formatting of
# the last data row down two more rows before overwriting values.
$ws.Range("A5:H5").Copy($ws.Range("A6:H6"))
$ws.Range("A5:H5").Copy($ws.Range("A7:H7"))

# Many of the numeric-looking columns (fund code, AUM, position %, etc.)
# are stored as plain text in the source data (to keep leading zeros /
# trailing zeros intact), so force those ranges to Text before writing.
$ws.Range("B2:G7").NumberFormat = "@"

# Row 2: 233009 / 大摩多因子精选策略混合
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "233009"
$ws.Range("C2").Value = "大摩多因子精选策略混合"
$ws.Range("D2").Value = "6.50"
$ws.Range("E2").Value = "83.44"
$ws.Range("F2").Value = "1.77"
$ws.Range("G2").Value = "0.1150"
$ws.Range("H2").Value = 1

# Row 3: 160613 / 鹏华盛世创新混合（LOF）
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "160613"
$ws.Range("C3").Value = "鹏华盛世创新混合（LOF）"
$ws.Range("D3").Value = "2.53"
$ws.Range("E3").Value = "92.08"
$ws.Range("F3").Value = "4.43"
$ws.Range("G3").Value = "0.1121"
$ws.Range("H3").Value = 6

# Row 4: 011574 / 鹏华领航一年持有期混合A
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "011574"
$ws.Range("C4").Value = "鹏华领航一年持有期混合A"
$ws.Range("D4").Value = "1.20"
$ws.Range("E4").Value = "92.84"
$ws.Range("F4").Value = "5.22"
$ws.Range("G4").Value = "0.0626"
$ws.Range("H4").Value = 5

# Row 5: 011575 / 鹏华领航一年持有期混合C
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "011575"
$ws.Range("C5").Value = "鹏华领航一年持有期混合C"
$ws.Range("D5").Value = "0.91"
$ws.Range("E5").Value = "92.84"
$ws.Range("F5").Value = "5.22"
$ws.Range("G5").Value = "0.0475"
$ws.Range("H5").Value = 5

# Row 6: 009384 / 摩根士丹利华鑫MSCI中国A股指数增强A
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "009384"
$ws.Range("C6").Value = "摩根士丹利华鑫MSCI中国A股指数增强A"
$ws.Range("D6").Value = "0.39"
$ws.Range("E6").Value = "90.98"
$ws.Range("F6").Value = "1.09"
$ws.Range("G6").Value = "0.0043"
$ws.Range("H6").Value = 6

# Row 7: 014866 / 摩根士丹利华鑫MSCI中国A股指数增强C
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "014866"
$ws.Range("C7").Value = "摩根士丹利华鑫MSCI中国A股指数增强C"
$ws.Range("D7").Value = "0.00"
$ws.Range("E7").Value = "90.98"
$ws.Range("F7").Value = "1.09"
# G7 is stored as a plain number (0), unlike the other G-column cells.
$ws.Range("G7").NumberFormat = "General"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 6

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet
# ---------------------------------------------------------------------
# Push the existing data rows (2-7) down to (3-8) - a straight block copy
# preserves each row's original formatting (including column A's bold
# bordered style). Then clone row 3's formatting back onto the freed-up
# row 2 (A2 from A3's style, B2:D2 from B3:D3's style) before overwriting
# every cell with its final value below.
$summary.Range("A2:D7").Copy($summary.Range("A3:D8"))
$summary.Range("A3").Copy($summary.Range("A2"))
$summary.Range("B3:D3").Copy($summary.Range("B2:D2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.34

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 4
$summary.Range("D3").Value = 0.33

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 7
$summary.Range("D4").Value = 0.42

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 3
$summary.Range("D5").Value = 0.34

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 4
$summary.Range("D6").Value = 0.34

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 7
$summary.Range("D7").Value = 0.25

$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 1
$summary.Range("D8").Value = 0.05

Write-Output "2022-Q3 sheet added and 总计 updated"
